# Append a new paragraph at the very end of the document body, matching
# the formatting of the paragraph that already ends with "(1 dia)".
$d = $word.ActiveDocument

$paras = $d.Paragraphs
$lastPara = $paras.Last

# Move to the end of the document content and insert a new paragraph mark.
$endRange = $lastPara.Range
$endRange.Collapse(0)  # wdCollapseEnd
$endRange.InsertParagraphAfter()

# The freshly created paragraph is now the last paragraph in the document.
$newPara = $d.Paragraphs.Last
$newRange = $newPara.Range

# Match formatting: Arial 12pt (sz 24 half-points), justified, 1.5 line spacing.
$newRange.Font.Name = "Arial"
$newRange.Font.Size = 12
$newPara.Format.Alignment = 3       # wdAlignParagraphJustify
$newPara.Format.LineSpacingRule = 1 # wdLineSpace1pt5

# Set the text of the new (empty) paragraph.
$newRange.Text = "Dia 16/09: 2hr (1 dia)"
